$d = $word.ActiveDocument

# --- Change 1: merge the two runs " " and "An technischen Ressourcen..." ---
# in the Budget paragraph into a single run (same visible text). We locate
# the paragraph by its distinctive text and rewrite its runs directly so
# only the trailing two runs are combined, leaving the first sentence run
# untouched.
$budgetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Auch das Zeitbudget*") {
        $budgetPara = $para
        break
    }
}

$budgetXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Auch das Zeitbudget ist anhand der Arbeitspakete aufgeteilt und in einem Google Tabellensheet festgehalten. Auf der linken Seite ist das Soll aufgeführt, auf der rechten Seite wird die tatsächlich investierte Zeit laufend festgehalten.</w:t></w:r><w:r><w:t xml:space="preserve"> An technischen Ressourcen steht für das Projekt ein Server der FHGR </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>mit Apache</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 2.4.25, auf dem PHP 7.4 läuft, zur Verfügung. </w:t></w:r></w:p>
'@

$budgetPara.Range.InsertXML($budgetXml)

# --- Change 2: rework the last paragraph ("Link zum GitLab-Board: ...") and
# append the new "Versionsverwaltung" section after it. ---
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:pPr><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">Link zum </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>GitLab</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">-Board: </w:t></w:r><w:hyperlink r:id="rId13"><w:r><w:rPr><w:color w:val="1155CC"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:u w:val="single"/></w:rPr><w:t>https://gitlab.com/yxaw/front-projekt/-/boards</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>$</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Versionsverwaltung</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Wir arbeiten mit den vier </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Branches</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>", "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>develop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>", "abgabe-prototyp" und "abgabe-projekt". Auf dem Branch "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>develop</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>" wird entwickelt,</w:t></w:r><w:r><w:t xml:space="preserve"> er enthält das </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>work</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-in-progress. F</w:t></w:r><w:r><w:t>unktionierende Zustände und Zwischenergebnisse werden auf den "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>" Branch gepusht. Der "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">" enthält </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ausserdem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> einen Ordner "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>documents</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>", welche</w:t></w:r><w:r><w:t>r</w:t></w:r><w:r><w:t xml:space="preserve"> zugehörige Dokumentationen beinhaltet. Die beiden anderen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Branches</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dienen den entsprechenden Abgaben.</w:t></w:r></w:p>
'@

$rng.InsertXML($xml)
